# Update "想去人数" (column F) counts on both the "展览" sheet and the
# aggregated "全部类型" sheet to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# (row, newValue) pairs for worksheet "展览"
$exhibitionUpdates = @(
    @{Row = 3;  Value = 56}
    @{Row = 5;  Value = 183}
    @{Row = 6;  Value = 9558}
    @{Row = 7;  Value = 858}
    @{Row = 8;  Value = 332}
    @{Row = 9;  Value = 1208}
    @{Row = 10; Value = 1698}
    @{Row = 11; Value = 157}
    @{Row = 12; Value = 106}
    @{Row = 14; Value = 271}
    @{Row = 15; Value = 452}
    @{Row = 18; Value = 1321}
)

# (row, newValue) pairs for worksheet "全部类型"
$allTypesUpdates = @(
    @{Row = 3;  Value = 56}
    @{Row = 6;  Value = 183}
    @{Row = 7;  Value = 9558}
    @{Row = 8;  Value = 858}
    @{Row = 9;  Value = 332}
    @{Row = 10; Value = 1208}
    @{Row = 11; Value = 1698}
    @{Row = 12; Value = 157}
    @{Row = 13; Value = 106}
    @{Row = 15; Value = 271}
    @{Row = 16; Value = 452}
    @{Row = 19; Value = 1321}
)

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($u in $exhibitionUpdates) {
    $wsExhibition.Cells.Item($u.Row, 6).Value = $u.Value
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($u in $allTypesUpdates) {
    $wsAllTypes.Cells.Item($u.Row, 6).Value = $u.Value
}
